$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (second worksheet)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status text update (affects every cell sharing this string, e.g. C2/C3 here
# and B/C columns on the Overview sheet).
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

# Latest Handback DateTime now populated (was the zero-date placeholder).
$wsZh.Range("H2").Value = "2016-03-24 17:19:59"
$wsZh.Range("H3").Value = "2016-03-24 17:19:59"

# New "Latest Target File" (F) / "Latest Handback File" (G) columns mirroring
# the handoff file info already present in columns A and D.
$wsZh.Range("F2").Value = "713d65fe-609f-47e7-a923-6cb5e0f3b62f.md"
$wsZh.Range("G2").Value = "713d65fe-609f-47e7-a923-6cb5e0f3b62f.bc8cfab82bf8d37015be60e53b56a676978fe7b4.zh-cn.xlf"
$wsZh.Range("F3").Value = "e6e3ebef-5325-48ef-8b2f-2d60623f08ab.md"
$wsZh.Range("G3").Value = "e6e3ebef-5325-48ef-8b2f-2d60623f08ab.681f5d9bea457651a691870722442bf7be47b0c7.zh-cn.xlf"

# Rebuild the hyperlinks collection in final document order so relationship
# ids come out as rId2..rId9 (rId1 stays the table part).
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0a7a54a40024dd2b4c1919592220ba64add4cf84/e2e/713d65fe-609f-47e7-a923-6cb5e0f3b62f.md", [Type]::Missing, [Type]::Missing, "713d65fe-609f-47e7-a923-6cb5e0f3b62f.md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d42ca057651b7e9911cee8bec793732c0ec5cac7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/713d65fe-609f-47e7-a923-6cb5e0f3b62f.bc8cfab82bf8d37015be60e53b56a676978fe7b4.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "713d65fe-609f-47e7-a923-6cb5e0f3b62f.bc8cfab82bf8d37015be60e53b56a676978fe7b4.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/0a7a54a40024dd2b4c1919592220ba64add4cf84/e2e/713d65fe-609f-47e7-a923-6cb5e0f3b62f.md", [Type]::Missing, [Type]::Missing, "713d65fe-609f-47e7-a923-6cb5e0f3b62f.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d42ca057651b7e9911cee8bec793732c0ec5cac7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/713d65fe-609f-47e7-a923-6cb5e0f3b62f.bc8cfab82bf8d37015be60e53b56a676978fe7b4.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "713d65fe-609f-47e7-a923-6cb5e0f3b62f.bc8cfab82bf8d37015be60e53b56a676978fe7b4.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0a7a54a40024dd2b4c1919592220ba64add4cf84/e2e/e6e3ebef-5325-48ef-8b2f-2d60623f08ab.md", [Type]::Missing, [Type]::Missing, "e6e3ebef-5325-48ef-8b2f-2d60623f08ab.md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d42ca057651b7e9911cee8bec793732c0ec5cac7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e6e3ebef-5325-48ef-8b2f-2d60623f08ab.681f5d9bea457651a691870722442bf7be47b0c7.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "e6e3ebef-5325-48ef-8b2f-2d60623f08ab.681f5d9bea457651a691870722442bf7be47b0c7.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/0a7a54a40024dd2b4c1919592220ba64add4cf84/e2e/e6e3ebef-5325-48ef-8b2f-2d60623f08ab.md", [Type]::Missing, [Type]::Missing, "e6e3ebef-5325-48ef-8b2f-2d60623f08ab.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d42ca057651b7e9911cee8bec793732c0ec5cac7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e6e3ebef-5325-48ef-8b2f-2d60623f08ab.681f5d9bea457651a691870722442bf7be47b0c7.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "e6e3ebef-5325-48ef-8b2f-2d60623f08ab.681f5d9bea457651a691870722442bf7be47b0c7.zh-cn.xlf")

# New cells get the same visual style (underlined hyperlink font) used by the
# rest of the hyperlink columns on this sheet.
$wsZh.Range("F2:G3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet "de-de" (third worksheet)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Range("H2").Value = "2016-03-24 17:20:10"
$wsDe.Range("H3").Value = "2016-03-24 17:20:10"

$wsDe.Range("F2").Value = "713d65fe-609f-47e7-a923-6cb5e0f3b62f.md"
$wsDe.Range("G2").Value = "713d65fe-609f-47e7-a923-6cb5e0f3b62f.bc8cfab82bf8d37015be60e53b56a676978fe7b4.de-de.xlf"
$wsDe.Range("F3").Value = "e6e3ebef-5325-48ef-8b2f-2d60623f08ab.md"
$wsDe.Range("G3").Value = "e6e3ebef-5325-48ef-8b2f-2d60623f08ab.681f5d9bea457651a691870722442bf7be47b0c7.de-de.xlf"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0a7a54a40024dd2b4c1919592220ba64add4cf84/e2e/713d65fe-609f-47e7-a923-6cb5e0f3b62f.md", [Type]::Missing, [Type]::Missing, "713d65fe-609f-47e7-a923-6cb5e0f3b62f.md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3599e40643f782ab31b51057381d440150324f90/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/713d65fe-609f-47e7-a923-6cb5e0f3b62f.bc8cfab82bf8d37015be60e53b56a676978fe7b4.de-de.xlf", [Type]::Missing, [Type]::Missing, "713d65fe-609f-47e7-a923-6cb5e0f3b62f.bc8cfab82bf8d37015be60e53b56a676978fe7b4.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/0a7a54a40024dd2b4c1919592220ba64add4cf84/e2e/713d65fe-609f-47e7-a923-6cb5e0f3b62f.md", [Type]::Missing, [Type]::Missing, "713d65fe-609f-47e7-a923-6cb5e0f3b62f.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3599e40643f782ab31b51057381d440150324f90/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/713d65fe-609f-47e7-a923-6cb5e0f3b62f.bc8cfab82bf8d37015be60e53b56a676978fe7b4.de-de.xlf", [Type]::Missing, [Type]::Missing, "713d65fe-609f-47e7-a923-6cb5e0f3b62f.bc8cfab82bf8d37015be60e53b56a676978fe7b4.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0a7a54a40024dd2b4c1919592220ba64add4cf84/e2e/e6e3ebef-5325-48ef-8b2f-2d60623f08ab.md", [Type]::Missing, [Type]::Missing, "e6e3ebef-5325-48ef-8b2f-2d60623f08ab.md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3599e40643f782ab31b51057381d440150324f90/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e6e3ebef-5325-48ef-8b2f-2d60623f08ab.681f5d9bea457651a691870722442bf7be47b0c7.de-de.xlf", [Type]::Missing, [Type]::Missing, "e6e3ebef-5325-48ef-8b2f-2d60623f08ab.681f5d9bea457651a691870722442bf7be47b0c7.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/0a7a54a40024dd2b4c1919592220ba64add4cf84/e2e/e6e3ebef-5325-48ef-8b2f-2d60623f08ab.md", [Type]::Missing, [Type]::Missing, "e6e3ebef-5325-48ef-8b2f-2d60623f08ab.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3599e40643f782ab31b51057381d440150324f90/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e6e3ebef-5325-48ef-8b2f-2d60623f08ab.681f5d9bea457651a691870722442bf7be47b0c7.de-de.xlf", [Type]::Missing, [Type]::Missing, "e6e3ebef-5325-48ef-8b2f-2d60623f08ab.681f5d9bea457651a691870722442bf7be47b0c7.de-de.xlf")

$wsDe.Range("F2:G3").Style = "HyperLink"
